# Generate Report for Handback
# Fills in the "Latest Target File", "Latest Handback File",
# "Latest Handback DateTime" and "Error Detail" columns for the row
# describing 748277fe-64a8-43d2-a9cf-f8784faea75e.md on both the
# zh-cn and de-de language sheets, now that the handback for that file
# has actually come in (but turned out to be stale / not the latest
# version), and widens the Error Detail column so the message is
# readable.

$wb = $excel.ActiveWorkbook

$targetMdName = "748277fe-64a8-43d2-a9cf-f8784faea75e.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2fef75a2af545b802273a04596effc0315ea7fa/e2e/748277fe-64a8-43d2-a9cf-f8784faea75e.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2d6eed091125f2f7f0bcdfa011461c1588d6551/e2e/748277fe-64a8-43d2-a9cf-f8784faea75e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2fef75a2af545b802273a04596effc0315ea7fa/e2e/748277fe-64a8-43d2-a9cf-f8784faea75e.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I5").Value = $targetMdName
$wsZh.Range("I5").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("I5"), $targetUrl, "", "", $targetMdName)

$wsZh.Range("J5").Value = "748277fe-64a8-43d2-a9cf-f8784faea75e.7371a624d50c21d67a5464e934be36f77fb6c49c.zh-cn.xlf"
$wsZh.Range("K5").Value = "2016-08-31 12:28:39"
$wsZh.Range("P5").Value = $errorDetail

# ColumnWidth 39.17 round-trips to the OOXML stored "width" of 40 (same
# character-width <-> stored-width conversion already used by every other
# width="40" column in this workbook).
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I5").Value = $targetMdName
$wsDe.Range("I5").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("I5"), $targetUrl, "", "", $targetMdName)

$wsDe.Range("J5").Value = "748277fe-64a8-43d2-a9cf-f8784faea75e.7371a624d50c21d67a5464e934be36f77fb6c49c.de-de.xlf"
$wsDe.Range("K5").Value = "2016-08-31 12:28:57"
$wsDe.Range("P5").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = 39.17
